# chore: update Sheets via scheduled runner
# Refreshes the cached market-board figures (currentAveragePrice* /
# LevePrice* / LeveProfit*, columns H:N) for the affected leve rows across
# all eight crafting-job tabs (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1000
$ws.Range("I19").Value = 1000
$ws.Range("J19").Value = 1000
$ws.Range("K19").Value = 1000
$ws.Range("L19").Value = 1000
$ws.Range("M19").Value = -825
$ws.Range("N19").Value = -1350
$ws.Range("H96").Value = 356.66666
$ws.Range("I96").Value = 265.66666
$ws.Range("K96").Value = 796.9999799999999
$ws.Range("M96").Value = 576.0000200000001
$ws.Range("H111").Value = 422.14285
$ws.Range("I111").Value = 475.83334
$ws.Range("J111").Value = 100
$ws.Range("K111").Value = 1427.50002
$ws.Range("L111").Value = 300
$ws.Range("M111").Value = 1639.49998
$ws.Range("N111").Value = -6434
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").Value = ""
$ws.Range("H116").Value = 4499.3335
$ws.Range("I116").Value = 4499.3335
$ws.Range("K116").Value = 4499.3335
$ws.Range("M116").Value = -1057.3335
$ws.Range("H118").Value = 0
$ws.Range("I118").Value = 0
$ws.Range("K118").Value = 0
$ws.Range("M118").Value = ""
$ws.Range("H138").Value = 3783.0312
$ws.Range("I138").Value = 2840
$ws.Range("J138").Value = 3957.6667
$ws.Range("K138").Value = 8520
$ws.Range("L138").Value = 11873.0001
$ws.Range("M138").Value = -3380
$ws.Range("N138").Value = -22153.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1421
$ws.Range("I2").Value = 1338.2858
$ws.Range("K2").Value = 1338.2858
$ws.Range("M2").Value = -1225.2858
$ws.Range("H116").Value = 1421
$ws.Range("I116").Value = 1338.2858
$ws.Range("K116").Value = 1338.2858
$ws.Range("M116").Value = 955.7141999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1421
$ws.Range("I3").Value = 1338.2858
$ws.Range("K3").Value = 1338.2858
$ws.Range("M3").Value = -1224.2858
$ws.Range("H86").Value = 5330
$ws.Range("I86").Value = 2827.6667
$ws.Range("J86").Value = 7832.3335
$ws.Range("K86").Value = 2827.6667
$ws.Range("L86").Value = 7832.3335
$ws.Range("M86").Value = -1704.6667
$ws.Range("N86").Value = -10078.3335
$ws.Range("H89").Value = 5330
$ws.Range("I89").Value = 2827.6667
$ws.Range("J89").Value = 7832.3335
$ws.Range("K89").Value = 14138.3335
$ws.Range("L89").Value = 39161.6675
$ws.Range("M89").Value = -8522.333500000001
$ws.Range("N89").Value = -50393.6675
$ws.Range("H94").Value = 611.8125
$ws.Range("I94").Value = 557.2143
$ws.Range("K94").Value = 557.2143
$ws.Range("M94").Value = -106.2143
$ws.Range("H99").Value = 995.75
$ws.Range("I99").Value = 996.3333
$ws.Range("J99").Value = 994
$ws.Range("K99").Value = 996.3333
$ws.Range("L99").Value = 994
$ws.Range("M99").Value = 501.6667
$ws.Range("N99").Value = -3990
$ws.Range("H105").Value = 2516
$ws.Range("I105").Value = 2518
$ws.Range("K105").Value = 2518
$ws.Range("M105").Value = -771
$ws.Range("H134").Value = 2182.5652
$ws.Range("I134").Value = 2066.611
$ws.Range("K134").Value = 6199.833
$ws.Range("M134").Value = -3664.833

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = ""
$ws.Range("N31").Value = ""
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = ""
$ws.Range("N34").Value = ""
$ws.Range("H58").Value = 8698
$ws.Range("I58").Value = 2396
$ws.Range("K58").Value = 2396
$ws.Range("M58").Value = -2193
$ws.Range("H136").Value = 8698
$ws.Range("I136").Value = 2396
$ws.Range("K136").Value = 7188
$ws.Range("M136").Value = -4638
$ws.Range("H141").Value = 49710.5
$ws.Range("J141").Value = 59947.332
$ws.Range("L141").Value = 59947.332
$ws.Range("N141").Value = -70307.33199999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H24").Value = 498
$ws.Range("I24").Value = 498
$ws.Range("K24").Value = 1494
$ws.Range("M24").Value = -1264
$ws.Range("H34").Value = 4342.3335
$ws.Range("J34").Value = 4342.3335
$ws.Range("L34").Value = 13027.0005
$ws.Range("N34").Value = -13195.0005
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").Value = ""
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").Value = ""
$ws.Range("H98").Value = 2291.25
$ws.Range("J98").Value = 3191.5
$ws.Range("L98").Value = 9574.5
$ws.Range("N98").Value = -12570.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").Value = ""
$ws.Range("H119").Value = 120000
$ws.Range("J119").Value = 120000
$ws.Range("L119").Value = 120000
$ws.Range("N119").Value = -129676

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H19").Value = 6669999.5
$ws.Range("I19").Value = 6669999.5
$ws.Range("K19").Value = 6669999.5
$ws.Range("M19").Value = -6669829.5
$ws.Range("H22").Value = 7254.5
$ws.Range("I22").Value = 5347.5
$ws.Range("J22").Value = 7731.25
$ws.Range("K22").Value = 5347.5
$ws.Range("L22").Value = 7731.25
$ws.Range("M22").Value = -5052.5
$ws.Range("N22").Value = -8321.25
$ws.Range("H27").Value = 7254.5
$ws.Range("I27").Value = 5347.5
$ws.Range("J27").Value = 7731.25
$ws.Range("K27").Value = 5347.5
$ws.Range("L27").Value = 7731.25
$ws.Range("M27").Value = -5240.5
$ws.Range("N27").Value = -7945.25
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = ""
$ws.Range("N40").Value = ""
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").Value = ""

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 512.1667
$ws.Range("I100").Value = 407.2
$ws.Range("J100").Value = 1037
$ws.Range("K100").Value = 814.4
$ws.Range("L100").Value = 2074
$ws.Range("M100").Value = -273.4
$ws.Range("N100").Value = -3156
$ws.Range("H101").Value = 45000
$ws.Range("J101").Value = 45000
$ws.Range("L101").Value = 45000
$ws.Range("N101").Value = -51490
